$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(100).Insert()

$ws.Cells.Item(100, 1).Value = 10
$ws.Cells.Item(100, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(100, 3).Value = "La Araucanía"
$ws.Cells.Item(100, 4).Value = 45128
$ws.Cells.Item(100, 5).Value = 9
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100107
$ws.Cells.Item(100, 8).Value = "Otros"
$ws.Cells.Item(100, 9).Value = 100107002
$ws.Cells.Item(100, 10).Value = "Chirimoya"
$ws.Cells.Item(100, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(100, 12).Value = "Primera"
$ws.Cells.Item(100, 13).Value = 25
$ws.Cells.Item(100, 14).Value = 5000
$ws.Cells.Item(100, 15).Value = 5000
$ws.Cells.Item(100, 16).Value = 5000
$ws.Cells.Item(100, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(100, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(100, 19).Value = 5000
$ws.Cells.Item(100, 20).Value = 1
